$d = $word.ActiveDocument

# 1. The second paragraph (currently empty) gets a new sentence of text.
$p2 = $d.Paragraphs(2)
$p2.Range.InsertAfter("InsDelF → deletions longer than 1 base are almost exclusively 2 base deletions, some from tandem seq, some w/ microhomology AG/CT and XXXX")

# 2. Insert a run of new paragraphs right after it (before the "Here, for the
#    first time..." paragraph). InsertParagraphAfter() clones the pPr of the
#    paragraph it is called on, which matches the Normal/480-line-spacing
#    paragraphs used throughout this section.
$newParaTexts = @(
    "InsDel4 → deletions of 2, 3, and 4 bases – mostly tandem repeats, not so many involving AG/CT",
    "microhom e.g. AAA → A/TTT→T",
    "TAT → T",
    "ATA → A",
    "AGGAG → AG",
    "AGAAG → AG",
    "",
    "",
    ""
)

$anchor = $p2
foreach ($text in $newParaTexts) {
    $anchor.Range.InsertParagraphAfter()
    $anchor = $d.Paragraphs($anchor.Index + 1)
    if ($text -ne "") {
        $anchor.Range.InsertAfter($text)
    }
}

# 3. Trim the leading sentence off what is now the paragraph following the
#    newly-inserted block (formerly paragraph 3, the "Here, for the first
#    time..." paragraph), leaving the rest of its content untouched.
$target = $d.Paragraphs($anchor.Index + 1)
$lead = "Here, for the first time, we identified C_ID4 and ID_F, along with their corresponding 89-type representations (InsDel4a, InsDel4b, and InsDel_F), using a de novo extraction approach. "
$start = $target.Range.Start
$leadRange = $d.Range($start, $start + $lead.Length)
$leadRange.Delete()
